$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 3.755628166162433

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 109.9114832445916
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 127.7687514715588

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5.582307763322248

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 6.15379541431027

$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 3.223369029078222
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 8.656069925401464

$ws.Range("B7").Value = 0.6545652718822623
$ws.Range("C7").Value = 0.3048912486333797
$ws.Range("D7").Value = 0.7210945179870265
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2.213936997104367

$ws.Range("B8").Value = 0.003078177322033415
$ws.Range("C8").Value = 0.04103571897497393
$ws.Range("D8").Value = 3.223369029078222
$ws.Range("E8").Value = 0.5333859586016987
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.800868883976928

$ws.Range("B9").Value = 0.6545652718822623
$ws.Range("C9").Value = 0.3048912486333797
$ws.Range("D9").Value = 0.7210945179870265
$ws.Range("E9").Value = 0.5333859586016987
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.213936997104367

$ws.Range("B10").Value = 1.445647641019636
$ws.Range("C10").Value = 1.626987699542094
$ws.Range("D10").Value = 3.223369029078222
$ws.Range("E10").Value = 0.5333859586016987
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.82939032824165

$ws.Range("B11").Value = 1.445647641019636
$ws.Range("C11").Value = 1.626987699542094
$ws.Range("D11").Value = 0.1496068669990043
$ws.Range("E11").Value = 0.5333859586016987
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.755628166162433

$ws.Range("B12").Value = 0.6545652718822623
$ws.Range("C12").Value = 0.3048912486333797
$ws.Range("D12").Value = 0.7210945179870265
$ws.Range("E12").Value = 0.5333859586016987
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 2.213936997104367

$ws.Range("B13").Value = 0.6545652718822623
$ws.Range("C13").Value = 1.626987699542094
$ws.Range("D13").Value = 0.7210945179870265
$ws.Range("E13").Value = 0.5333859586016987
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 3.536033448013082

$ws.Range("B14").Value = 3.272327238179451
$ws.Range("C14").Value = 1.626987699542094
$ws.Range("D14").Value = 3.223369029078222
$ws.Range("E14").Value = 13.86384647080068
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 21.98653043760045

$ws.Range("B15").Value = 3.272327238179451
$ws.Range("C15").Value = 1.626987699542094
$ws.Range("D15").Value = 0.7210945179870265
$ws.Range("E15").Value = 0.5333859586016987
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 6.15379541431027

$ws.Range("B16").Value = 3.272327238179451
$ws.Range("C16").Value = 1.626987699542094
$ws.Range("D16").Value = 18.71679738969934
$ws.Range("E16").Value = 0.5333859586016987
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 24.14949828602258

$ws.Range("B17").Value = 3.272327238179451
$ws.Range("C17").Value = 1.626987699542094
$ws.Range("D17").Value = 0.7210945179870265
$ws.Range("E17").Value = 0.5333859586016987
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 6.15379541431027

$ws.Range("B18").Value = 3.272327238179451
$ws.Range("C18").Value = 1.626987699542094
$ws.Range("D18").Value = 0.1496068669990043
$ws.Range("E18").Value = 0.5333859586016987
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 5.582307763322248

$ws.Range("B19").Value = 3.272327238179451
$ws.Range("C19").Value = 1.626987699542094
$ws.Range("D19").Value = 0.7210945179870265
$ws.Range("E19").Value = 0.5333859586016987
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 6.15379541431027
